$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the data range to Text format before writing, so numeric-looking
# strings (prices like "323.32") are stored as text, matching the source
# workbook (all data cells are inlineStr/shared-string text, not numbers).
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "27.355.78"
$ws.Range("E2").Value = "  -3.66%  "
$ws.Range("D3").Value = "1.857.95"
$ws.Range("E3").Value = "  -4.62%  "
$ws.Range("E4").Value = "  -1.15%  "
$ws.Range("D5").Value = "323.32"
$ws.Range("E5").Value = "  +0.54%  "
$ws.Range("E6").Value = "  -0.98%  "
$ws.Range("D7").Value = "0.4529"
$ws.Range("E7").Value = "  -5.53%  "
$ws.Range("D8").Value = "0.3862"
$ws.Range("E8").Value = "  -5.56%  "
$ws.Range("D9").Value = "48.58"
$ws.Range("E9").Value = "  -9.38%  "
$ws.Range("D10").Value = "0.07911"
$ws.Range("E10").Value = "  -6.97%  "
$ws.Range("E11").Value = "  -3.69%  "
$ws.Range("D12").Value = "21.39"
$ws.Range("E12").Value = "  -4.33%  "
$ws.Range("D13").Value = "1.868.49"
$ws.Range("E13").Value = "  -4.20%  "
$ws.Range("D14").Value = "5.904"
$ws.Range("E14").Value = "  -4.01%  "
$ws.Range("D15").Value = "7.122"
$ws.Range("E15").Value = "  -5.79%  "
$ws.Range("E16").Value = "  -1.19%  "
$ws.Range("D17").Value = "0.00001032"
$ws.Range("E17").Value = "  -3.81%  "
$ws.Range("D18").Value = "85.77"
$ws.Range("E18").Value = "  -4.90%  "
$ws.Range("E19").Value = "  -1.93%  "
$ws.Range("D20").Value = "17.10"
$ws.Range("E20").Value = "  -6.96%  "
$ws.Range("E21").Value = "  -0.98%  "
$ws.Range("D22").Value = "5.529"
$ws.Range("E22").Value = "  -5.24%  "
$ws.Range("D23").Value = "27.362.39"
$ws.Range("E23").Value = "  -3.77%  "
$ws.Range("D24").Value = "10.84"
$ws.Range("D25").Value = "2.281"
$ws.Range("E25").Value = "  -0.76%  "
$ws.Range("D26").Value = "2.087.44"
$ws.Range("E26").Value = "  -4.47%  "
$ws.Range("D27").Value = "153.91"
$ws.Range("E27").Value = "  -1.43%  "
$ws.Range("D28").Value = "19.74"
$ws.Range("E28").Value = "  -2.72%  "
$ws.Range("D29").Value = "2.068"
$ws.Range("E29").Value = "  -4.98%  "
$ws.Range("D30").Value = "5.436"
$ws.Range("E30").Value = "  -6.74%  "
$ws.Range("D31").Value = "120.65"
$ws.Range("E31").Value = "  -2.82%  "
$ws.Range("D32").Value = "1.486"
$ws.Range("E32").Value = "  +3.48%  "
$ws.Range("D33").Value = "0.09280"
$ws.Range("E33").Value = "  -3.91%  "
$ws.Range("D34").Value = "0.9339"
$ws.Range("E34").Value = "  -5.08%  "
$ws.Range("D35").Value = "3.595"
$ws.Range("E35").Value = "  -2.85%  "
$ws.Range("D36").Value = "5.258"
$ws.Range("E36").Value = "  -6.40%  "
$ws.Range("E37").Value = "  -4.03%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").Value = "1.221"
$ws.Range("E38").Value = "  -2.09%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "0.05988"
$ws.Range("E39").Value = "  -3.15%  "
$ws.Range("D40").Value = "8.213"
$ws.Range("E40").Value = "  -10.41%  "
$ws.Range("E41").Value = "  -0.95%  "
$ws.Range("D42").Value = "0.5905"
$ws.Range("E42").Value = "  -5.07%  "
$ws.Range("D43").Value = "0.1887"
$ws.Range("E43").Value = "  -1.72%  "
$ws.Range("E44").Value = "  -9.73%  "
$ws.Range("D45").Value = "1.277"
$ws.Range("E45").Value = "  -4.07%  "
$ws.Range("D46").Value = "0.5617"
$ws.Range("E46").Value = "  -5.58%  "
$ws.Range("D47").Value = "11.94"
$ws.Range("E47").Value = "  -7.72%  "
$ws.Range("D48").Value = "3.373"
$ws.Range("E48").Value = "  -1.01%  "
$ws.Range("E49").Value = "  -6.63%  "
$ws.Range("D50").Value = "0.06773"
$ws.Range("E50").Value = "  -0.52%  "
$ws.Range("D51").Value = "108.10"
$ws.Range("E51").Value = "  -2.05%  "

# Remove the temporary Text number-format again so the cells end up with
# no explicit style index, same as in the source file.
$dataRange.ClearFormats()
